$d = $word.ActiveDocument

# The document (before edit) has 5 paragraphs:
#   1. Title "Share Files between Windows and Ubuntu"
#   2. (empty)
#   3. "Command: sudo mount -a -t vboxsf Share '/home/shayan/Share'"
#   4. (empty)
#   5. "Notice: Once the command is entered ... in both directions."
#      (this paragraph also carries the trailing _GoBack bookmark markers)
#
# The edit removes the trailing empty paragraph and the whole "Notice" text,
# while keeping the _GoBack bookmark markers, which end up re-attached to
# the end of the "Command" paragraph (paragraph 3).

# Step 1: remove the empty paragraph that sits right before "Notice".
# At this point it is paragraph 4 (not yet the document's last paragraph),
# so deleting its Range cleanly removes that paragraph mark.
$emptyPara = $d.Paragraphs.Item(4)
$emptyPara.Range.Delete()

# Step 2: clear the "Notice: ..." sentence text (now paragraph 4, and the
# document's last paragraph). Using wildcard Find/Replace removes the runs'
# text while leaving the paragraph mark (and the bookmark inside it) intact.
$d.Content.Find.Execute("Notice:*directions.", $true, $false, $true, $false, $false, $true, 1, $false, "", 2) | Out-Null

# Step 3: merge the now-empty last paragraph into the "Command" paragraph
# (paragraph 3) by deleting paragraph 3's own trailing paragraph mark. This
# pulls the bookmark (which lives in the now-empty final paragraph) up to
# immediately follow the "/Share'" run, matching the target document.
$cmdPara = $d.Paragraphs.Item(3)
$markRange = $d.Range($cmdPara.Range.End - 1, $cmdPara.Range.End)
$markRange.Delete()
